# "added timely care to MHBC"
#
# Row 2 holds the parameters used to build the mod_* calls in rows 4-6:
#   A2 = Unique Resource Name, C2 = Name of Data subset (others stay the same).
# Rename the PMH-based entry to the BCMH-based "TimelyCare" entry:
#   A2: PMH_TimelyCare -> BCMH_TimelyCare
#   C2: PMHdata()      -> MHBC()
# Rows 4-6 are formulas referencing A2/B2/C2/D2/E2/F2/G2, so their cached
# values update automatically once row 2 changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BCMH_TimelyCare"
$ws.Range("C2").Value = "MHBC()"

# Rows 10-12 used to hold stale, hard-coded copies of the old row 4-6
# results. Turn them into live formulas (mirroring rows 4-6, but keyed off
# row 8, which is currently blank) so they behave like the rest of the sheet.
$ws.Range("B10").Formula = '="mod_Accordion_ui("&"''"&A8&"''"&")"'
$ws.Range("B11").Formula = '="mod_Accordion_server("&"''"&A8&"'', selector=selection, data="&C8&", title = c(''"&D8&"''), Visible = T)"'
$ws.Range("B12").Formula = '="mod_info_server(''"&A8&"'', selector = selection, data = "&C8&", rownametitle = c(''"&B8&"''), phone = "&E8&", website = "&F8&", email = "&G8&")"'

# Move the active selection to B21, matching where the author ended up.
$ws.Range("B21").Select() | Out-Null
